# Update Sage scrape results
#
# 1. Insert a new column before column G ("The" / Estonia 2007 shifts to H,
#    etc. all the way through to Sandowrm landing in Z). The new data-row
#    cells in the inserted column are filled with 0 to match the other
#    "hit count" columns.
# 2. Refresh the bibliographic details (Title / Year / DOI / Access Type)
#    for every record row with the latest scrape results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "The" column before G -------------------------------
$ws.Columns("G:G").Insert()
$ws.Range("G1").Value = "The"
$ws.Range("G2:G11").Value = 0

# Helper: write a value into a cell as plain TEXT, even when the text looks
# like a pure number (e.g. "2015"), without leaving a lasting style/format
# change on the cell (matches the cell's original un-styled inlineStr type).
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# --- 2. Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value = "Technical note: exploiting problem definition study for cyber security simulations"
Set-TextValue "D2" "2015"
Set-TextValue "E2" "10.1177/1548512915604585"

# --- Row 3 -------------------------------------------------------------------
$ws.Range("B3").Value = "The dynamics of cyber conflict between rival antagonists, 2001–11"
Set-TextValue "D3" "2014"
Set-TextValue "E3" "10.1177/0022343313518940"

# --- Row 4 -------------------------------------------------------------------
$ws.Range("B4").Value = "A virtual necessity: Some modest steps toward greater cybersecurity"
Set-TextValue "D4" "2016"
Set-TextValue "E4" "10.1177/0096340212459039"

# --- Row 5 -------------------------------------------------------------------
$ws.Range("B5").Value = "Ontological security, cyber technology, and states’ responses"
Set-TextValue "D5" "2022"
Set-TextValue "E5" "10.1177/13540661221130958"
$ws.Range("F5").Value = "Open Access"

# --- Row 6 -------------------------------------------------------------------
$ws.Range("B6").Value = "Cyclones in cyberspace: Information shaping and denial in the 2008 Russia–Georgia war"
Set-TextValue "D6" "2012"
Set-TextValue "E6" "10.1177/0967010611431079"

# --- Row 7 -------------------------------------------------------------------
$ws.Range("B7").Value = "On 3D simultaneous attack against manoeuvring target with communication delays"
Set-TextValue "D7" "2020"
Set-TextValue "E7" "10.1177/1729881419894808"
$ws.Range("F7").Value = "Restricted"

# --- Row 8 -------------------------------------------------------------------
$ws.Range("B8").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
Set-TextValue "D8" "2021"
Set-TextValue "E8" "10.1177/10439862211001606"
$ws.Range("F8").Value = "Restricted"

# --- Row 9 -------------------------------------------------------------------
$ws.Range("B9").Value = "Moving beyond the sanctuary paradigm: Canada must face up to the reality of a contested and dangerous space environment"
Set-TextValue "D9" "2023"
Set-TextValue "E9" "10.1177/00207020231178394"

# --- Row 10 ------------------------------------------------------------------
$ws.Range("B10").Value = "The code not taken: China, the United States, and the future of cyber espionage"
Set-TextValue "D10" "2013"
Set-TextValue "E10" "10.1177/0096340213501344"
$ws.Range("F10").Value = "Restricted"

# --- Row 11 ------------------------------------------------------------------
$ws.Range("B11").Value = "Using network digital twins to improve cyber resilience of missions"
Set-TextValue "D11" "2022"
Set-TextValue "E11" "10.1177/15485129221131226"
